$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "302.02"
Set-TextCell "E2" "-0.67%"
Set-TextCell "D3" "37.43"
Set-TextCell "E3" "7.02%"
Set-TextCell "D4" "4.995"
Set-TextCell "E4" "-3.29%"
Set-TextCell "D5" "0.07822"
Set-TextCell "E5" "0.38%"
Set-TextCell "D6" "2.192"
Set-TextCell "E6" "-8.69%"
Set-TextCell "D7" "8.035"
Set-TextCell "D8" "4.036"
Set-TextCell "E8" "1.62%"
Set-TextCell "D9" "0.9092"
Set-TextCell "E9" "-2.73%"
Set-TextCell "D10" "0.09697"
Set-TextCell "E10" "-2.98%"
Set-TextCell "D11" "0.1894"
Set-TextCell "E11" "2.54%"
Set-TextCell "D12" "0.08506"
Set-TextCell "E12" "-2.23%"
Set-TextCell "D13" "0.03526"
Set-TextCell "E13" "6.20%"
Set-TextCell "D14" "0.09971"
Set-TextCell "E14" "0.79%"
Set-TextCell "D15" "0.001490"
Set-TextCell "E15" "-0.33%"
Set-TextCell "D16" "0.005650"
Set-TextCell "E16" "-1.63%"
Set-TextCell "D17" "3.469"
Set-TextCell "E17" "0.05%"
Set-TextCell "E18" "-3.70%"
Set-TextCell "E19" "2.54%"
Set-TextCell "D20" "0.1293"
Set-TextCell "E20" "-2.57%"
Set-TextCell "D21" "4.767"
Set-TextCell "E22" "-0.95%"
Set-TextCell "D23" "0.04636"
Set-TextCell "D24" "0.001230"
Set-TextCell "E24" "1.05%"
Set-TextCell "D25" "0.004803"
Set-TextCell "E25" "8.27%"
Set-TextCell "E26" "-7.57%"
Set-TextCell "D27" "0.0004753"
Set-TextCell "E27" "28.56%"
Set-TextCell "D39" "0.01757"
Set-TextCell "E39" "-1.28%"
Set-TextCell "D40" "0.04725"
Set-TextCell "E40" "-1.82%"
Set-TextCell "D41" "0.008060"
Set-TextCell "E41" "3.97%"
Set-TextCell "D42" "0.1394"
Set-TextCell "E42" "-1.22%"
Set-TextCell "D43" "0.007672"
Set-TextCell "E43" "7.65%"
Set-TextCell "D44" "0.002171"
Set-TextCell "E44" "-0.70%"
Set-TextCell "D45" "0.01040"
Set-TextCell "E45" "13.37%"
Set-TextCell "D46" "0.00006057"
Set-TextCell "E46" "2.19%"
Set-TextCell "D47" "0.00000000751"
Set-TextCell "E47" "0.15%"
Set-TextCell "D48" "6.066"
Set-TextCell "E48" "121.99%"
Set-TextCell "D50" "0.00002101"
Set-TextCell "E50" "0.15%"
Set-TextCell "D51" "0.0002001"
Set-TextCell "E51" "0.15%"
